$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3691694736480713
$ws.Range("B1").Value = 0.2613443732261658
$ws.Range("C1").Value = 0.4093986749649048
$ws.Range("D1").Value = 4.563246726989746
$ws.Range("E1").Value = 2.29588770866394
